# Insert a new weekly data row for "Apio" (Feria Lagunitas de Puerto Montt)
# right above the existing row 454. This pushes the existing rows 454-531
# down to 455-532 (dimension grows from A1:R531 to A1:R532) and the new
# row 454 carries the newest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 454:531 down by one row, leaving row 454 blank for the new entry.
$ws.Rows("454:454").Insert()

$newRow = 454

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 45180
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112017
$ws.Cells.Item($newRow, 7).Value = "Apio"
$ws.Cells.Item($newRow, 8).Value = "Americana (o)"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 15
$ws.Cells.Item($newRow, 11).Value = 12000
$ws.Cells.Item($newRow, 12).Value = 12000
$ws.Cells.Item($newRow, 13).Value = 12000
$ws.Cells.Item($newRow, 14).Value = "`$/docena de matas"
$ws.Cells.Item($newRow, 15).Value = "Región de Coquimbo"
$ws.Cells.Item($newRow, 16).Value = 2000
$ws.Cells.Item($newRow, 17).Value = 6
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
